$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.456.55"
$ws.Range("E2").Value = "  -3.91%  "
$ws.Range("D3").Value = "3.028.61"
$ws.Range("E3").Value = "  -6.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.67"
$ws.Range("E5").Value = "  -5.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.03"
$ws.Range("E6").Value = "  -9.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.010.44"
$ws.Range("E8").Value = "  -6.94%  "
$ws.Range("E9").Value = "  -12.25%  "
$ws.Range("E10").Value = "  -11.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.99"
$ws.Range("E11").Value = "  -10.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -11.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "34.19"
$ws.Range("E13").Value = "  -13.39%  "
$ws.Range("E14").Value = "  -14.13%  "
$ws.Range("D15").Value = "3.530.28"
$ws.Range("E15").Value = "  -6.11%  "
$ws.Range("D16").Value = "64.488.19"
$ws.Range("E16").Value = "  -3.93%  "
$ws.Range("E17").Value = "  -3.98%  "
$ws.Range("D18").Value = "3.034.72"
$ws.Range("E18").Value = "  -6.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "483.16"
$ws.Range("E19").Value = "  -9.71%  "
$ws.Range("E20").Value = "  -10.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.39"
$ws.Range("E21").Value = "  -11.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.658"
$ws.Range("E22").Value = "  -14.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.94"
$ws.Range("E23").Value = "  -12.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.39"
$ws.Range("E24").Value = "  -11.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "77.74"
$ws.Range("E25").Value = "  -9.89%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -15.04%  "
$ws.Range("E28").Value = "  -7.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.46"
$ws.Range("E29").Value = "  -9.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.61"
$ws.Range("E30").Value = "  -13.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.54"
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.12"
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "507.88"
$ws.Range("E34").Value = "  -5.74%  "
$ws.Range("E35").Value = "  -9.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.74"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.76"
$ws.Range("E37").Value = "  -13.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0404"
$ws.Range("E38").Value = "  -5.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0794"
$ws.Range("E39").Value = "  -9.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.118"
$ws.Range("E40").Value = "  -8.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.11"
$ws.Range("E41").Value = "  -13.99%  "
$ws.Range("D42").Value = "2.805.72"
$ws.Range("E42").Value = "  -5.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.37"
$ws.Range("E43").Value = "  -11.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.237"
$ws.Range("E45").Value = "  -12.14%  "
$ws.Range("E46").Value = "  -8.03%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.106"
$ws.Range("E47").Value = "  -7.85%  "
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").Value = "0.0₃0509"
$ws.Range("E48").Value = "  -14.78%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.03"
$ws.Range("E49").Value = "  -5.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.42"
$ws.Range("E50").Value = "  -12.44%  "
$ws.Range("E51").Value = "  -18.04%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
